{"js": "// Xu\u1ea5t kho b\u00e1n tr\u1ef1u ti\u1ebfp / xem tr\u01b0\u1edbc\n//\n// 1. Drop the old `_GoBack` bookmark that sat right after the\n//    `loaiHinhKho` merge field.\n// 2. Clear the (both/justify) paragraph alignment on the \"M\u00e3 c\u00e2n\" and\n//    \"S\u1ed1 bao b\u00ec\" sample-row cells.\n// 3. Right-align the \"Tr\u1ecdng l\u01b0\u1ee3ng c\u1ea3 b\u00ec\" sample-row cell.\n// 4. Re-drop `_GoBack` at the spot the cursor ended up \u2014 in the middle\n//    of the `trongLuongCaBi` merge-field placeholder text \u2014 which\n//    splits that run in two.\n\n// --- 1. remove the stale _GoBack bookmark -------------------------------\ncontext.document.deleteBookmark(\"_GoBack\");\nawait context.sync();\n\n// --- 2/3. fix up alignment on the sample data row of the 2nd table -----\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst dataTable = tables.items[1];\nconst rows = dataTable.rows;\nrows.load(\"items\");\nawait context.sync();\n\nconst dataRow = rows.items[1];\nconst cells = dataRow.cells;\ncells.load(\"items\");\nawait context.sync();\n\nconst maCanParagraph = cells.items[1].body.paragraphs;\nconst soBaoBiParagraph = cells.items[2].body.paragraphs;\nconst trongLuongCaBiParagraph = cells.items[3].body.paragraphs;\nmaCanParagraph.load(\"items\");\nsoBaoBiParagraph.load(\"items\");\ntrongLuongCaBiParagraph.load(\"items\");\nawait context.sync();\n\n// \"Left\" serializes as \"no <w:jc> at all\" (left is the implicit default),\n// matching the diff which just drops the element.\nmaCanParagraph.items[0].alignment = \"Left\";\nsoBaoBiParagraph.items[0].alignment = \"Left\";\ntrongLuongCaBiParagraph.items[0].alignment = \"Right\";\nawait context.sync();\n\n// --- 4. split the trongLuongCaBi placeholder run & re-insert _GoBack ---\nconst firstHalf = context.document.body.search(\n  \"\u00ab#if($!d.trongLuongCaBi)$numbe\",\n  { matchCase: true }\n);\nfirstHalf.load(\"items\");\nawait context.sync();\n\nconst splitPoint = firstHalf.items[0].getRange(\"End\");\nsplitPoint.insertBookmark(\"_GoBack\");\nawait context.sync();\n", "ps1": "# Xu\u1ea5t kho b\u00e1n tr\u1ef1u ti\u1ebfp / xem tr\u01b0\u1edbc\n#\n# 1. Drop the old `_GoBack` bookmark that sat right after the\n#    `loaiHinhKho` merge field.\n# 2. Clear the (both/justify) paragraph alignment on the \"M\u00e3 c\u00e2n\" and\n#    \"S\u1ed1 bao b\u00ec\" sample-row cells.\n# 3. Right-align the \"Tr\u1ecdng l\u01b0\u1ee3ng c\u1ea3 b\u00ec\" sample-row cell.\n# 4. Re-drop `_GoBack` at the spot the cursor ended up \u2014 in the middle\n#    of the `trongLuongCaBi` merge-field placeholder text \u2014 which\n#    splits that run in two.\n\n$d = $word.ActiveDocument\n\n# --- 1. remove the stale _GoBack bookmark -------------------------------\nif ($d.Bookmarks.Exists(\"_GoBack\")) {\n    $d.Bookmarks.Item(\"_GoBack\").Delete()\n}\n\n# --- 2/3. fix up alignment on the sample data row of the 2nd table -----\n$dataTable = $d.Tables.Item(2)\n$dataRow = $dataTable.Rows.Item(2)\n\n# wdAlignParagraphLeft = 0 -> serializes with no <w:jc> at all (left is\n# the implicit default), matching the diff which just drops the element.\n$dataRow.Cells.Item(2).Range.Paragraphs.Item(1).Alignment = 0\n$dataRow.Cells.Item(3).Range.Paragraphs.Item(1).Alignment = 0\n# wdAlignParagraphRight = 2\n$dataRow.Cells.Item(4).Range.Paragraphs.Item(1).Alignment = 2\n\n# --- 4. split the trongLuongCaBi placeholder run & re-insert _GoBack ---\n$findRange = $d.Content\n$findRange.Find.ClearFormatting()\n$needle = [char]0x00AB + \"#if(`$!d.trongLuongCaBi)`$numbe\"\n$found = $findRange.Find.Execute($needle)\nif ($found) {\n    $findRange.Collapse(0)  # wdCollapseEnd\n    $d.Bookmarks.Add(\"_GoBack\", $findRange)\n}\n"}
